$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.186.77'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.781.81'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.87'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.70'
$ws.Range('E8').Value = '  -1.01%  '
$ws.Range('E9').Value = '  +0.40%  '
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0946'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.037.72'
$ws.Range('E12').Value = '  +0.40%  '
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.775.17'
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.128.02'
$ws.Range('E15').Value = '  +1.25%  '
$ws.Range('E16').Value = '  +2.17%  '
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.85'
$ws.Range('E18').Value = '  +2.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.96'
$ws.Range('E19').Value = '  +3.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0789'
$ws.Range('E20').Value = '  +1.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.97'
$ws.Range('E21').Value = '  +3.71%  '
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('E23').Value = '  +2.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.04'
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.95'
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.17'
$ws.Range('E26').Value = '  +2.21%  '
$ws.Range('E27').Value = '  +1.35%  '
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.01'
$ws.Range('E29').Value = '  +0.47%  '
$ws.Range('E30').Value = '  +0.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0518'
$ws.Range('E31').Value = '  +1.41%  '
$ws.Range('E32').Value = '  +3.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.64'
$ws.Range('E33').Value = '  +4.13%  '
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.442.39'
$ws.Range('E35').Value = '  +4.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.651'
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('E37').Value = '  +8.05%  '
$ws.Range('E38').Value = '  +3.78%  '
$ws.Range('E39').Value = '  +1.39%  '
$ws.Range('E40').Value = '  +0.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '80.16'
$ws.Range('E41').Value = '  +2.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.923'
$ws.Range('E42').Value = '  +1.78%  '
$ws.Range('E43').Value = '  +0.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.59'
$ws.Range('E44').Value = '  +0.45%  '
$ws.Range('E45').Value = '  +2.09%  '
$ws.Range('E46').Value = '  +3.69%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.07'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0137'
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.940.41'
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '104.75'
$ws.Range('E50').Value = '  -2.24%  '
$ws.Range('E51').Value = '  +0.27%  '
